$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-19 12:03:24"
$wsOverview.Range("G3").Value = "2016-10-19 12:03:24"

# --- zh-cn sheet ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-19 12:03:12"
$wsZhCn.Range("H3").Value = "2016-10-19 12:03:12"
$wsZhCn.Range("K2").Value = "2016-10-19 12:03:54"
$wsZhCn.Range("K3").Value = "2016-10-19 12:03:54"

# --- de-de sheet ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-19 12:03:24"
$wsDeDe.Range("H3").Value = "2016-10-19 12:03:24"
$wsDeDe.Range("K2").Value = "2016-10-19 12:04:13"
$wsDeDe.Range("K3").Value = "2016-10-19 12:04:13"
